$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.354.98'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.874.18'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'0.7142"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").Value = "'241.83"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.3110"
$ws.Range("E8").Value = '  +1.02%  '
$ws.Range("D9").Value = "'0.07759"
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = "'25.11"
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '1.878.36'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = "'5.256"
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("D14").Value = "'0.7126"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = "'91.16"
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '29.363.39'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = "'6.084"
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").Value = "'0.000008226"
$ws.Range("E18").Value = '  +5.19%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").Value = '2.124.71'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = "'7.788"
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = "'0.1595"
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D27").Value = "'9.050"
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").Value = "'18.53"
$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").Value = "'4.418"
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").Value = "'4.332"
$ws.Range("E31").Value = '  +2.74%  '
$ws.Range("D32").Value = "'1.283"
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("D33").Value = "'0.05311"
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").Value = "'1.938"
$ws.Range("E34").Value = '  +1.41%  '
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("D36").Value = "'0.7401"
$ws.Range("E36").Value = '  -9.12%  '
$ws.Range("D37").Value = "'2.699"
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").Value = "'0.01870"
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("D39").Value = '1.232.35'
$ws.Range("E39").Value = '  +5.24%  '
$ws.Range("D40").Value = "'2.730"
$ws.Range("E40").Value = '  +1.13%  '
$ws.Range("D41").Value = "'6.521"
$ws.Range("E41").Value = '  +5.49%  '
$ws.Range("D42").Value = "'110.50"
$ws.Range("E42").Value = '  +8.49%  '
$ws.Range("D43").Value = "'0.8914"
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("D44").Value = "'73.01"
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '2.023.29'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = "'1.809"
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").Value = "'9.453"
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("D51").Value = "'0.4315"
$ws.Range("E51").Value = '  +1.23%  '
